$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 40-41: Bittensor and FirstDigitalUSD swap positions
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"

$ws.Range("D2").Value = "59.578.10"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "2.609.73"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.566"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.51"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.335"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").Value = "3.068.65"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("D14").Value = "59.501.90"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000134"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "2.559.13"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "340.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.409"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.166"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.71%  "
$ws.Range("D28").Value = "0.0₃0745"
$ws.Range("E28").Value = "  +3.17%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +5.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("E35").Value = "  +0.79%  "
$ws.Range("E36").Value = "  +3.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.826"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "277.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.602"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0949"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("D46").Value = "1.940.68"
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.80%  "
